# Update occupation classification code
# - insert a new "% of total people" column (E), pushing "average debt per
#   person" to column F
# - reorder the occupation rows so "merchant" is first
# - populate the new "% of total people" values

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Insert a new column before E, shifting "average debt per person" (and
#    its data) to F. Range.EntireColumn.Insert shifts cells right and carries
#    formatting from the column at the insertion point, matching the header
#    cell style (bold/centered/bordered).
$ws.Range("E1").EntireColumn.Insert()

# 2) Header row
$ws.Range("E1").Value = "% of total people"

# 3) Capture the current (pre-reorder) data rows 2-7 so we can rewrite them
#    in the new order: merchant moves to the top, the rest keep their
#    relative order.
$rows = @{}
for ($r = 2; $r -le 7; $r++) {
    $rows[$r] = @{
        B = $ws.Cells.Item($r, 2).Value()
        C = $ws.Cells.Item($r, 3).Value()
        D = $ws.Cells.Item($r, 4).Value()
        F = $ws.Cells.Item($r, 6).Value()
    }
}

# Old layout (row => occupation): 2 broker, 3 doctor, 4 farmer, 5 mariner,
# 6 merchant, 7 schoolmaster.
# New layout: 2 merchant, 3 broker, 4 doctor, 5 farmer, 6 mariner,
# 7 schoolmaster.
$newOrder = @(6, 2, 3, 4, 5, 7)

$percents = @{
    2 = 80
    3 = 3.333333333333333
    4 = 3.333333333333333
    5 = 6.666666666666667
    6 = 3.333333333333333
    7 = 3.333333333333333
}

$destRow = 2
foreach ($srcRow in $newOrder) {
    $data = $rows[$srcRow]
    # Column A is just the sequential row index (0-5) and is unaffected by
    # the reorder.
    $ws.Cells.Item($destRow, 2).Value = $data.B
    $ws.Cells.Item($destRow, 3).Value = $data.C
    $ws.Cells.Item($destRow, 4).Value = $data.D
    $ws.Cells.Item($destRow, 5).Value = $percents[$destRow]
    $ws.Cells.Item($destRow, 6).Value = $data.F
    $destRow++
}
